# This workbook is an Aspose.Cells "Evaluation Warning" output sample.
# The only meaningful content change in this revision is the refreshed
# copyright year stamped into the warning sheet's text (2014 -> 2016),
# which comes from regenerating the example with a newer Aspose.Cells
# build. Update that shared string via the "Evaluation Warning" sheet.
$wb = $excel.ActiveWorkbook

$warningSheet = $wb.Worksheets.Item("Evaluation Warning")
$warningSheet.Range("A5").Value = "Evaluation Only. Created with Aspose.Cells for .NET.Copyright 2003 - 2016 Aspose Pty Ltd."
